# Insert a new data row at row 111 (pushing existing rows 111-176 down to 112-177)
# and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(111).Insert()

$ws.Cells.Item(111, 1).Value = 10
$ws.Cells.Item(111, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(111, 3).Value = "La Araucanía"
$ws.Cells.Item(111, 4).Value = 44755
$ws.Cells.Item(111, 5).Value = 9
$ws.Cells.Item(111, 6).Value = 100112013
$ws.Cells.Item(111, 7).Value = "Alcachofa"
$ws.Cells.Item(111, 8).Value = "Madrigal"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 200
$ws.Cells.Item(111, 11).Value = 18000
$ws.Cells.Item(111, 12).Value = 18000
$ws.Cells.Item(111, 13).Value = 18000
$ws.Cells.Item(111, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(111, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(111, 16).Value = 450
$ws.Cells.Item(111, 17).Value = 40
$ws.Cells.Item(111, 18).Value = "Hortaliza"
